# Applies the edit described by the diff.
#
# The diff consists of two real text edits:
#   1. "... của Bộ Xây dựng của các phương tiện nêu tại Điều 1 ..."
#      -> "... của Bộ Xây dựng đối với các phương tiện nêu tại Điều 1 ..."
#   2. "- Công an tỉnh;                                  (p/h)"
#      -> "- Phòng CSGT Công an tỉnh;                                  (p/h)"
#
# ...plus a number of purely cosmetic clean-ups where a run that had been
# split in two/three pieces around a <w:proofErr/> (grammar-check) marker is
# consolidated back into a single run with no visible text change:
#   - "... tỉnh Bắc Ninh;" (end of the "Căn cứ Quyết định ..." clause)
#   - "Vận tải & An toàn giao thông" (three separate places)
#   - "- Như Điều 3;"
#   - "- Sở Xây dựng các tỉnh, thành phố;"
#   - "Lý do thu hồi"
#   - "{ghi_chu}" / "{/don_vi_list}"
#
# Doing a Find/Replace over the full run of text (even replacing text with
# itself) makes Word rebuild the run(s) it touches as a single run and drops
# any <w:proofErr/> markers inside the replaced range, which reproduces the
# consolidation seen in the diff.

$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindContinue = 1

# 1) "Căn cứ Quyết định ..." clause: merge trailing "Bắc Ninh;" back in.
$d.Content.Find.Execute(
    "Căn cứ Quyết định số 12/2025/QĐ-UBND ngày 01/7/2025 của Uỷ ban nhân dân tỉnh Bắc Ninh ban hành Quy định về chức năng, nhiệm vụ, quyền hạn của Sở Xây dựng tỉnh Bắc Ninh;",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "Căn cứ Quyết định số 12/2025/QĐ-UBND ngày 01/7/2025 của Uỷ ban nhân dân tỉnh Bắc Ninh ban hành Quy định về chức năng, nhiệm vụ, quyền hạn của Sở Xây dựng tỉnh Bắc Ninh;",
    $wdReplaceAll)

# 2) "Theo đề nghị của Trưởng phòng Vận tải & An toàn giao thông." (italic run)
$d.Content.Find.Execute(
    "Theo đề nghị của Trưởng phòng Vận tải & An toàn giao thông.",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "Theo đề nghị của Trưởng phòng Vận tải & An toàn giao thông.",
    $wdReplaceAll)

# 3) "Giao Phòng Vận tải & An toàn giao thông ..." + "của" -> "đối với"
$d.Content.Find.Execute(
    "Giao Phòng Vận tải & An toàn giao thông thực hiện chuyển trạng thái phù hiệu",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "Giao Phòng Vận tải & An toàn giao thông thực hiện chuyển trạng thái phù hiệu",
    $wdReplaceAll)

$d.Content.Find.Execute(
    "của Bộ Xây dựng của các phương tiện nêu tại Điều 1",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "của Bộ Xây dựng đối với các phương tiện nêu tại Điều 1",
    $wdReplaceAll)

# 4) "Chánh Văn phòng Sở, Trưởng phòng Vận tải & An toàn giao thông, ..."
$d.Content.Find.Execute(
    "Chánh Văn phòng Sở, Trưởng phòng Vận tải & An toàn giao thông,",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "Chánh Văn phòng Sở, Trưởng phòng Vận tải & An toàn giao thông,",
    $wdReplaceAll)

# 5) "- Như Điều 3;"
$d.Content.Find.Execute(
    "- Như Điều 3;",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "- Như Điều 3;",
    $wdReplaceAll)

# 6) "- Sở Xây dựng các tỉnh, thành phố;"
$d.Content.Find.Execute(
    "- Sở Xây dựng các tỉnh, thành phố;",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "- Sở Xây dựng các tỉnh, thành phố;",
    $wdReplaceAll)

# 7) "- Công an tỉnh; ... (p/h)" -> "- Phòng CSGT Công an tỉnh; ... (p/h)"
$d.Content.Find.Execute(
    "- Công an tỉnh;",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "- Phòng CSGT Công an tỉnh;",
    $wdReplaceAll)

# 8) "Lý do thu hồi" (table header)
$d.Content.Find.Execute(
    "Lý do thu hồi",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "Lý do thu hồi",
    $wdReplaceAll)

# 9) "{ghi_chu}" / "{/don_vi_list}" template placeholders (different runs/sizes)
$d.Content.Find.Execute(
    "{ghi_chu}",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "{ghi_chu}",
    $wdReplaceAll)

$d.Content.Find.Execute(
    "{/don_vi_list}",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "{/don_vi_list}",
    $wdReplaceAll)
